# Fruta / hortaliza, semanal
# A new weekly observation is inserted at row 51 (pushing the existing
# rows 51-107 down to 52-108), with the original data set now occupying
# rows 52-108 (i.e. each pre-existing row shifts down by one position).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 51; Excel shifts rows 51..107 down to 52..108
# and extends the sheet's used range/dimension automatically.
$ws.Rows.Item(51).Insert()

# Clone the formatting/content of the row directly above (row 50) into the
# newly inserted (now blank) row 51, so every constant column (Mercado,
# Region, Codreg, Tipo, Producto, Categoria, Variedad, Calidad, Unidad,
# Kg/unidad, etc.) keeps the same values/styles used throughout the table.
$ws.Rows.Item(50).Copy()
$ws.Rows.Item(51).PasteSpecial()

# Now overwrite the columns that hold the new weekly record's own data.
$ws.Cells.Item(51, 4).Value = 44671    # D51 Fecha
$ws.Cells.Item(51, 13).Value = 400     # M51 Volumen
$ws.Cells.Item(51, 14).Value = 7500    # N51 Precio minimo
$ws.Cells.Item(51, 15).Value = 8000    # O51 Precio maximo
$ws.Cells.Item(51, 16).Value = 7750    # P51 Precio promedio ponderado
$ws.Cells.Item(51, 18).Value = "Ecuador" # R51 Origen
$ws.Cells.Item(51, 19).Value = 1938    # S51 Precio $/Kg
